$wb = $excel.ActiveWorkbook

# Sheet 1 -> "2025"
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 1037.265132737054
$ws.Range("E2").Value = 28926.05393052954
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 16171.06685703679
$ws.Range("L2").Value = 48492.22142001599
$ws.Range("M2").Value = 10595.37713982
$ws.Range("N2").Value = 7085.795531257033
$ws.Range("O2").Value = 6997.710127123046

# Sheet 2 -> "2030"
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 4157.588990853394
$ws.Range("E2").Value = 45991.90904307188
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 37079.12819938764
$ws.Range("L2").Value = 54844.03303316472
$ws.Range("M2").Value = 17449.04999683176
$ws.Range("N2").Value = 9043.260164362724
$ws.Range("O2").Value = 9733.027027076689

# Sheet 3 -> "2035"
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 2754.31755456332
$ws.Range("B2").Value = 6368.910634126893
$ws.Range("E2").Value = 57457.45307013817
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 52465.73681402855
$ws.Range("L2").Value = 54844.03303316472
$ws.Range("M2").Value = 21912.87293902603
$ws.Range("N2").Value = 13060.5833893223
$ws.Range("O2").Value = 12869.82982044365

# Sheet 4 -> "2040"
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = 2754.31755456332
$ws.Range("B2").Value = 6368.910634126893
$ws.Range("E2").Value = 57457.45307013817
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 52465.73681402855
$ws.Range("L2").Value = 54844.03303316472
$ws.Range("M2").Value = 21912.87293902603
$ws.Range("N2").Value = 13178.75918555561
$ws.Range("O2").Value = 12869.82982044365

# Sheet 5 -> "2045"
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = 5713.151062849596
$ws.Range("B2").Value = 6368.910634126893
$ws.Range("E2").Value = 57457.45307013817
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 52465.73681402855
$ws.Range("L2").Value = 54844.03303316472
$ws.Range("M2").Value = 21912.87293902603
$ws.Range("N2").Value = 13630.65979676794
$ws.Range("O2").Value = 14953.0932418782

# Sheet 6 -> "2050"
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = 5713.151062849596
$ws.Range("B2").Value = 6368.910634126893
$ws.Range("E2").Value = 57457.45307013817
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 52465.73681402855
$ws.Range("L2").Value = 54844.03303316472
$ws.Range("M2").Value = 21912.87293902603
$ws.Range("N2").Value = 13630.65979676794
$ws.Range("O2").Value = 14953.0932418782
